# Auto-generated Excel COM-interop edit script
# Applies the cryptos-list price/volume refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (no numeric auto-detection / no lost trailing zeros)
# for the price cells whose new value still parses as a plain number.
foreach ($addr in @("D5","D6","D8","D17","D19","D20","D23","D25","D26","D27","D28","D31","D41","D43","D45","D48","D49")) {
    $ws.Range($addr).NumberFormat = "@"
}

# Cell value updates
$ws.Range("D2").Value = "27.129.28"
$ws.Range("E2").Value = "  +0.82%  "
$ws.Range("D3").Value = "1.567.61"
$ws.Range("E3").Value = "  +1.05%  "
$ws.Range("E4").Value = "  +0.80%  "
$ws.Range("D5").Value = "210.46"
$ws.Range("D6").Value = "0.491"
$ws.Range("E6").Value = "  +0.68%  "
$ws.Range("E7").Value = "  +0.75%  "
$ws.Range("D8").Value = "21.99"
$ws.Range("E8").Value = "  -0.20%  "
$ws.Range("E9").Value = "  +0.21%  "
$ws.Range("E10").Value = "  +0.47%  "
$ws.Range("E11").Value = "  +0.65%  "
$ws.Range("D12").Value = "1.788.92"
$ws.Range("E12").Value = "  +0.93%  "
$ws.Range("D13").Value = "1.567.76"
$ws.Range("E13").Value = "  +1.22%  "
$ws.Range("E14").Value = "  +0.45%  "
$ws.Range("E15").Value = "  +0.08%  "
$ws.Range("D16").Value = "27.113.18"
$ws.Range("E16").Value = "  +0.80%  "
$ws.Range("D17").Value = "62.01"
$ws.Range("E17").Value = "  +0.70%  "
$ws.Range("E18").Value = "  -0.94%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "7.38"
$ws.Range("E19").Value = "  +1.16%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "214.69"
$ws.Range("E20").Value = "  -1.15%  "
$ws.Range("E21").Value = "  +0.81%  "
$ws.Range("E22").Value = "  +1.28%  "
$ws.Range("D23").Value = "9.20"
$ws.Range("E23").Value = "  +0.25%  "
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").Value = "154.10"
$ws.Range("E25").Value = "  +0.33%  "
$ws.Range("D26").Value = "6.60"
$ws.Range("E26").Value = "  -0.62%  "
$ws.Range("D27").Value = "15.05"
$ws.Range("E27").Value = "  +0.51%  "
$ws.Range("D28").Value = "0.106"
$ws.Range("E28").Value = "  +1.28%  "
$ws.Range("E29").Value = "  +0.68%  "
$ws.Range("E30").Value = "  +5.30%  "
$ws.Range("D31").Value = "0.0473"
$ws.Range("E31").Value = "  +0.64%  "
$ws.Range("E32").Value = "  +0.81%  "
$ws.Range("E33").Value = "  +3.06%  "
$ws.Range("D34").Value = "1.430.70"
$ws.Range("E34").Value = "  +1.33%  "
$ws.Range("E35").Value = "  +11.65%  "
$ws.Range("E36").Value = "  +0.76%  "
$ws.Range("E37").Value = "  +2.11%  "
$ws.Range("E38").Value = "  +0.80%  "
$ws.Range("E39").Value = "  +0.82%  "
$ws.Range("E40").Value = "  +3.05%  "
$ws.Range("D41").Value = "0.808"
$ws.Range("E41").Value = "  +0.18%  "
$ws.Range("E42").Value = "  +0.74%  "
$ws.Range("D43").Value = "2.36"
$ws.Range("E43").Value = "  +2.16%  "
$ws.Range("E44").Value = "  +0.37%  "
$ws.Range("D45").Value = "64.54"
$ws.Range("E45").Value = "  +0.15%  "
$ws.Range("E46").Value = "  +0.81%  "
$ws.Range("D47").Value = "1.707.10"
$ws.Range("E47").Value = "  +1.22%  "
$ws.Range("D48").Value = "85.94"
$ws.Range("E48").Value = "  -1.34%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "0.0519"
$ws.Range("E49").Value = "  -0.40%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0100"
$ws.Range("E50").Value = "  -0.22%  "
$ws.Range("E51").Value = "  -0.11%  "
